$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.666.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.919.95"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.69"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4938"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2986"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06761"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.917.98"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.22"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07348"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.191"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.82"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6744"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.645.26"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007954"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.52"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.59%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.158.86"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.404"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +11.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "200.27"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.327"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.656"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.39"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.71"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.963"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.00%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.383"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09207"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.069"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05288"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7431"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.119"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.718"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.728"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9276"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.094"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4473"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.973"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.03"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +24.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.37"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.004"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1402"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.638"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.043"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.07"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.17%  "

# Rows 50 and 51 swap coin identity along with updated data
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05889"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.40%  "

$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4044"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.54%  "
